# Ticket 79/81/82 workbook update:
# Add a new worksheet "NoSpaceAfterParen" (used to test that JETT updates a
# cell reference in a formula even when no space precedes the opening
# parenthesis, e.g. IF(B2="-",0,B2) ).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet so it lands at the end
# of the tab order (matching "TagParseInFormula" -> "NoSpaceAfterParen").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "NoSpaceAfterParen"

# Header row: a / b / c / result - bold + centered.
$ws.Range("A1").Value = "a"
$ws.Range("B1").Value = "b"
$ws.Range("C1").Value = "c"
$ws.Range("D1").Value = "result"

$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Data / template row.
$ws.Range("A2").Value = '<jt:for start="1" end="10" var="x">${x}'
$ws.Range("B2").Value = '${x+1}'
$ws.Range("C2").Value = '${x+2}'
$ws.Range("D2").Value = '$[A2-(IF(B2="-",0,B2)+C2)]'
$ws.Range("E2").Value = '</jt:for>'
